# Daily attendance processing - 2025-12-12 19:04:46
# Reorders the semicolon-free, comma-separated list of "Recorded By" emails
# in column G for the affected rows (dedup / canonicalization of order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System"
    3  = "majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
    4  = "majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
    5  = "Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    6  = "manar.montaser@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
    7  = "menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
    12 = "yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, dina.adel@med.asu.edu.eg"
    13 = "yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
    17 = "mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg"
    19 = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
    24 = "youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
    27 = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
    28 = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
